$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-27 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-28 Sunday", 2)
$d.Content.Find.Execute("489×4=1956", $true, $false, $false, $false, $false, $true, 1, $false, "905×3=2715", 2)
$d.Content.Find.Execute("434×7=3038", $true, $false, $false, $false, $false, $true, 1, $false, "883×2=1766", 2)
$d.Content.Find.Execute("212×7=1484", $true, $false, $false, $false, $false, $true, 1, $false, "426×6=2556", 2)
$d.Content.Find.Execute("780×6=4680", $true, $false, $false, $false, $false, $true, 1, $false, "365×3=1095", 2)
$d.Content.Find.Execute("881×7=6167", $true, $false, $false, $false, $false, $true, 1, $false, "412×9=3708", 2)
$d.Content.Find.Execute("650×7=4550", $true, $false, $false, $false, $false, $true, 1, $false, "913×9=8217", 2)
$d.Content.Find.Execute("376×7=2632", $true, $false, $false, $false, $false, $true, 1, $false, "290×7=2030", 2)
$d.Content.Find.Execute("925×3=2775", $true, $false, $false, $false, $false, $true, 1, $false, "907×7=6349", 2)
$d.Content.Find.Execute("755×3=2265", $true, $false, $false, $false, $false, $true, 1, $false, "699×9=6291", 2)
$d.Content.Find.Execute("806×6=4836", $true, $false, $false, $false, $false, $true, 1, $false, "863×8=6904", 2)
$d.Content.Find.Execute("974×6=5844", $true, $false, $false, $false, $false, $true, 1, $false, "173×9=1557", 2)
$d.Content.Find.Execute("626×6=3756", $true, $false, $false, $false, $false, $true, 1, $false, "269×9=2421", 2)
$d.Content.Find.Execute("261×6=1566", $true, $false, $false, $false, $false, $true, 1, $false, "592×4=2368", 2)
$d.Content.Find.Execute("457×7=3199", $true, $false, $false, $false, $false, $true, 1, $false, "263×4=1052", 2)
$d.Content.Find.Execute("541×6=3246", $true, $false, $false, $false, $false, $true, 1, $false, "750×2=1500", 2)
$d.Content.Find.Execute("975×8=7800", $true, $false, $false, $false, $false, $true, 1, $false, "671×9=6039", 2)
$d.Content.Find.Execute("961×3=2883", $true, $false, $false, $false, $false, $true, 1, $false, "421×6=2526", 2)
$d.Content.Find.Execute("796×2=1592", $true, $false, $false, $false, $false, $true, 1, $false, "712×5=3560", 2)
$d.Content.Find.Execute("914×4=3656", $true, $false, $false, $false, $false, $true, 1, $false, "370×8=2960", 2)
$d.Content.Find.Execute("467×7=3269", $true, $false, $false, $false, $false, $true, 1, $false, "995×6=5970", 2)
$d.Content.Find.Execute("136×3=408", $true, $false, $false, $false, $false, $true, 1, $false, "512×3=1536", 2)
$d.Content.Find.Execute("693×9=6237", $true, $false, $false, $false, $false, $true, 1, $false, "330×3=990", 2)
$d.Content.Find.Execute("648×5=3240", $true, $false, $false, $false, $false, $true, 1, $false, "847×8=6776", 2)
$d.Content.Find.Execute("406×6=2436", $true, $false, $false, $false, $false, $true, 1, $false, "377×5=1885", 2)
$d.Content.Find.Execute("829×3=2487", $true, $false, $false, $false, $false, $true, 1, $false, "170×9=1530", 2)
